{"js": "// Exercise 2 (\"Teams pro Land\") gets its closing sentence reworded, and the\n// \"Derzeitige Laufzeit\" measurement right below it is updated from a quick\n// placeholder value to the real (much slower, unoptimised) measured time.\n\nconst body = context.document.body;\n\n// 1) Reword the tail of the \"Ein weiterer Diskussionspunkt ...\" paragraph\n//    (the paragraph that introduces Aufgabe 2 / \"Teams pro Land\").\nconst oldTail = \"eine Liste erstellen die L\u00e4nder und ihre Anzahl an Teams beinhalten soll.\";\nconst newTail =\n  \"gut vorbereitet sein und das Land mit den meisten Niederlassungen ausfindig machen. Die Abfrage jedoch ist nicht besonders hilfreich.\";\n\nconst tailHits = body.search(oldTail, { matchCase: true });\ntailHits.load(\"items\");\nawait context.sync();\n\nif (tailHits.items.length > 0) {\n  tailHits.items[0].insertText(newTail, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Update the \"Derzeitige Laufzeit:\" value that directly follows the\n//    Aufgabe-2 paragraph (the first \"9,6s\" in the document) to \"124,458 s\".\n//    The second \"9,6s\" later in the document (Aufgabe 3) must stay untouched.\nconst timeHits = body.search(\"9,6s\", { matchCase: true });\ntimeHits.load(\"items\");\nawait context.sync();\n\nif (timeHits.items.length > 0) {\n  timeHits.items[0].insertText(\"124,458 s\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Exercise 2 (\"Teams pro Land\") gets its closing sentence reworded, and the\n# \"Derzeitige Laufzeit\" measurement right below it is updated from a quick\n# placeholder value to the real (much slower, unoptimised) measured time.\n\n$d = $word.ActiveDocument\n\n# 1) Reword the tail of the \"Ein weiterer Diskussionspunkt ...\" paragraph\n#    (the paragraph that introduces Aufgabe 2 / \"Teams pro Land\").\n$oldTail = \"eine Liste erstellen die L\u00e4nder und ihre Anzahl an Teams beinhalten soll.\"\n$newTail = \"gut vorbereitet sein und das Land mit den meisten Niederlassungen ausfindig machen. Die Abfrage jedoch ist nicht besonders hilfreich.\"\n\n$find = $d.Content.Find\n$find.Text = $oldTail\n$find.Replacement.Text = $newTail\n# wdFindContinue=1 wrap, wdReplaceOne=1 -> replace only the single (unique) match\n$find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 1)\n\n# 2) Update the \"Derzeitige Laufzeit:\" value that directly follows the\n#    Aufgabe-2 paragraph (the first \"9,6s\" in the document) to \"124,458 s\".\n#    The second \"9,6s\" later in the document (Aufgabe 3) must stay untouched,\n#    so only a single replacement (wdReplaceOne) starting from the top of the\n#    document is performed.\n$oldTime = \"9,6s\"\n$newTime = \"124,458 s\"\n\n$find2 = $d.Content.Find\n$find2.Text = $oldTime\n$find2.Replacement.Text = $newTime\n$find2.Execute($oldTime, $true, $false, $false, $false, $false, $true, 1, $false, $newTime, 1)\n"}
